$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.455.83"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "2.109.84"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("E4").Value = "  +0.80%  "

$ws.Range("D5").Value = "'335.06"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("D7").Value = "'0.5228"
$ws.Range("E7").Value = "  +0.81%  "

$ws.Range("D8").Value = "'0.4552"
$ws.Range("E8").Value = "  +5.60%  "

$ws.Range("D9").Value = "'53.50"
$ws.Range("E9").Value = "  +17.04%  "

$ws.Range("D10").Value = "'0.08932"
$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("D11").Value = "'1.176"
$ws.Range("E11").Value = "  +2.08%  "

$ws.Range("D12").Value = "'24.32"
$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").Value = "2.105.98"
$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("D14").Value = "'6.857"
$ws.Range("E14").Value = "  +3.13%  "

$ws.Range("D15").Value = "'8.048"
$ws.Range("E15").Value = "  +5.06%  "

$ws.Range("D16").Value = "'96.72"
$ws.Range("E16").Value = "  +1.80%  "

$ws.Range("D17").Value = "'0.00001144"
$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").Value = "'1.007"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").Value = "'0.06654"
$ws.Range("E19").Value = "  +0.83%  "

$ws.Range("D20").Value = "'19.23"
$ws.Range("E20").Value = "  +2.76%  "

$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").Value = "'6.368"
$ws.Range("E22").Value = "  +2.61%  "

$ws.Range("D23").Value = "30.521.81"
$ws.Range("E23").Value = "  +0.81%  "

$ws.Range("D24").Value = "'12.42"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").Value = "'2.368"
$ws.Range("E25").Value = "  +4.06%  "

$ws.Range("D26").Value = "2.355.34"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("D27").Value = "'22.33"
$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("D28").Value = "'2.551"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").Value = "'163.54"
$ws.Range("E29").Value = "  +1.29%  "

$ws.Range("D30").Value = "'132.89"
$ws.Range("E30").Value = "  +1.72%  "

$ws.Range("D31").Value = "'1.227"
$ws.Range("E31").Value = "  +3.43%  "

$ws.Range("D32").Value = "'0.1072"
$ws.Range("E32").Value = "  +0.88%  "

$ws.Range("D33").Value = "'1.690"
$ws.Range("E33").Value = "  +8.77%  "

$ws.Range("D34").Value = "'6.345"
$ws.Range("E34").Value = "  +4.37%  "

$ws.Range("D35").Value = "'3.946"
$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("D36").Value = "'10.48"
$ws.Range("E36").Value = "  +9.18%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.725"
$ws.Range("E37").Value = "  +6.08%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02588"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").Value = "'0.06832"
$ws.Range("E39").Value = "  +3.55%  "

$ws.Range("D40").Value = "'0.2309"
$ws.Range("E40").Value = "  +3.50%  "

$ws.Range("D41").Value = "'12.73"
$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("D42").Value = "'0.6891"
$ws.Range("E42").Value = "  +3.19%  "

$ws.Range("D43").Value = "'1.249"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("D44").Value = "'2.338"
$ws.Range("E44").Value = "  +6.79%  "

$ws.Range("D45").Value = "'1.006"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").Value = "'14.13"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("D47").Value = "'0.6387"
$ws.Range("E47").Value = "  +1.28%  "

$ws.Range("D48").Value = "'3.668"
$ws.Range("E48").Value = "  +2.05%  "

$ws.Range("D49").Value = "'1.252"
$ws.Range("E49").Value = "  +1.59%  "

$ws.Range("D50").Value = "'0.3410"
$ws.Range("E50").Value = "  +25.40%  "

$ws.Range("D51").Value = "'83.31"
$ws.Range("E51").Value = "  +2.58%  "
